# Cmf.Calc.xlsx edit script
# - Rename the "Index" column header (table "testdata", column A) to "i"
# - Re-number the index column values from 1-based (1..502) to 0-based (0..501)
# - Narrow column A to fit the new, shorter content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of the first table column (was "Index", now "i").
# Setting the header cell's value on a ListObject/table column automatically
# renames the table column as well.
$ws.Range("A1").Value = "i"

# Shift the index column from 1-based to 0-based numbering for all data rows.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A no longer needs to fit the word "Index" (5 chars) or values up to
# 3 digits with a header driving the width; narrow it to match the new
# narrower content.
$ws.Columns.Item(1).ColumnWidth = 3.14
